# RS-RI schema doc update — align field descriptions with the new wording.
$d = $word.ActiveDocument

# Newline-within-cell marker (becomes a <w:br/> between <w:t> runs, same as
# existing manual line breaks in this document).
$br = [char]11

function Set-CellText($tableIndex, $row, $col, $text) {
    $table = $d.Tables.Item($tableIndex)
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
}

# --- Table 1 (caseId / mobilizedResource) ---------------------------------
Set-CellText 1 3 5 "Objet permettant de communquer la liste des ressource et vecteurs mobilisés"

# --- Table 2 (resource fields) ---------------------------------------------
Set-CellText 2 2 5 "A valoriser avec la date et heure d'engagement de la ressource/du vecteur"

$txt = "A valoriser avec l'identifiant partagé unique de la ressource engagée, normé comme suit :" + $br
$txt = $txt + "{orgID}.resource.{ID unique de la ressource partagée}" + $br
$txt = $txt + "OU - uniquement dans le cas où un ID unique de ressource ne peut pas être garanti par l'organisation propriétaire :" + $br
$txt = $txt + "{orgID}.resource.{sendercaseId}.{n° d’ordre chronologique de la ressource}"
Set-CellText 2 3 5 $txt

$txt = "A valoriser avec l'identifiant unique partagé de la demande de ressource (si la ressource a été engagée suite à une demande de ressource), normé comme suit :" + $br
$txt = $txt + "{orgID}.request.{ID unique de la demande dans le système émetteur}" + $br
$txt = $txt + "OU - si un ID unique de la demande n'était pas disponible : " + $br
$txt = $txt + "{OrgId émetteur}.request.{senderCaseId}.{numéro d’ordre chronologique}"
Set-CellText 2 4 5 $txt

Set-CellText 2 5 5 "A valoriser avec le numéro de mission unique du central d’appel (PSAP, …) qui a déclenché le vecteur"

$txt = "A valoriser avec l'identifiant de l'organisation à laquelle appartient la ressource, normé comme suit : " + $br
$txt = $txt + "{pays}.{domaine}.{organisation}"
Set-CellText 2 6 5 $txt

Set-CellText 2 7 5 "A valoriser avec le lieu de garage principal"

Set-CellText 2 8 5 "A valoriser avec le yype de ressource mobilisée : cf.nomenclature associée."

Set-CellText 2 9 5 "A valoriser avec le type de vecteur mobilisé : cf. nomenclature associée"

Set-CellText 2 10 5 "A valoriser avec le n° d'immatriculation du vecteur"

Set-CellText 2 11 5 "A valoriser avec le nom donné à la ressource par l’organisation d’appartenance"

Set-CellText 2 12 5 "A valoriser avec le code INSEE de la commune du centre d'affectation"

Set-CellText 2 13 5 "Objet qui décrit l'équipe à bord du vecteur"

Set-CellText 2 14 5 "Objet qui permet de décrire l'historique des états connu du vecteur mobilisé - et à minima le dernier état connu. "

Set-CellText 2 15 5 "Objet qui permet de transmettre la dernière géolocalisation connue d'un vecteur, au moment de la création du RS-RI."

Set-CellText 2 16 5 "A valoriser avec le type et valeur de l'URI utilisée par la ressource."

# --- Table 3 (team) ----------------------------------------------------
Set-CellText 3 2 5 "A valoriser avec le  niveau de médicalisation du vecteur. Cf. nomenclature associée"

Set-CellText 3 3 5 "A valoriser avec le nom de l'équipe à bord du vecteur (celui communiqué par l'organisation à laquelle l'équipe appartient)"

# --- Table 4 (state) -----------------------------------------------------
Set-CellText 4 2 5 "A valoriser avec la date et heure d'engagement de changement vers le nouveau statut"

Set-CellText 4 3 5 "A valoriser avec le statut du vecteur. Cf nomenclature associée."

$txt = "A valoriser de manière à indiquer la disponibilité du vecteur." + $br
$txt = $txt + "TRUE = DISPONIBLE" + $br
$txt = $txt + "FALSE = INDISPONIBLE" + $br
$txt = $txt + "VIDE = INCONNU"
Set-CellText 4 4 5 $txt

# --- Table 5 (coord) ------------------------------------------------------
Set-CellText 5 2 5 "A valoriser avec la latitude du point clé de la localisation "

Set-CellText 5 3 5 "A valoriser avec la longitude du point clé de la localisation"

Set-CellText 5 4 5 "A valoriser avec l'altitude du point clé de la localisation, en mètre, ignoré côté NexSIS. "

Set-CellText 5 5 5 "A valoriser en degrés"

Set-CellText 5 6 5 "A valoriser en km/h (notamment fournie par eCall, tel, nouveau AML)"

$txt = "A valoriser avec le niveau de précision des coordonnées fournies par le système emetteur. Cf. nomenclature associée." + $br
$txt = $txt + "CITY=Précision à l'échelle de la ville, STREET=Précision à l'échelle de la rue, ADDRESS=Adresse précise, EXACT=Point coordonnée GPS exact, UNKNOWN=Précision de la localisation non évaluable par l'émetteur"
Set-CellText 5 7 5 $txt
